# fixed player not participated in this season error
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Players in rows 11-14 had no value in the "Darts Used" (column G) field.
# Set it to "N/A" so the season stats don't error out for them.
$ws.Range("G11").Value = "N/A"
$ws.Range("G12").Value = "N/A"
$ws.Range("G13").Value = "N/A"
$ws.Range("G14").Value = "N/A"

# Update the sheet's active cell/selection.
$ws.Range("A16").Select()
